$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "Fandika Saputra" is replaced by a new hire "Susetyadi TP " ---
# New NIK + name, highlighted with a light-blue fill and black font.
$ws.Range("A8").Value = 4000224
$ws.Range("B8").Value = "Susetyadi TP "

$tmpStyle1 = $wb.Styles.Add("TmpHighlightStyle")
$tmpStyle1.Font.Color = 0
$tmpStyle1.Interior.Color = 16247773
$ws.Range("A8:B8").Style = "TmpHighlightStyle"
$wb.Styles.Item("TmpHighlightStyle").Delete()

# --- Row 10 (new row): another new hire "Brigitta d'Avriella" ---
# Same black font, but no fill this time.
$ws.Range("A10").Value = 4001344
$ws.Range("B10").Value = "Brigitta d'Avriella"

$tmpStyle2 = $wb.Styles.Add("TmpFontOnlyStyle")
$tmpStyle2.Font.Color = 0
$ws.Range("A10:B10").Style = "TmpFontOnlyStyle"
$wb.Styles.Item("TmpFontOnlyStyle").Delete()

# --- Move the active selection, matching the saved view state ---
$ws.Range("B13").Select() | Out-Null
